# Update the LR-pairs data to reflect the new TPM values, including a
# new "ECs" sending-cluster block (rows 2-5), shifting the previous data
# down and updating every numeric column throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Angptl3","Itgav","ECs",1,0.3333333333333333,0.030331,0.090993,0.005348220947889333,0.005348220947889334,3,1,9.423852333333334,28.271557,0.06654336290212845,0.06654336290212845,0.2858348651223334,2.572513786101,0.0003558886074161653,0.0003558886074161653),
    @("ECs","Angptl3","Itgav","FAPs",1,0.3333333333333333,0.030331,0.090993,0.005348220947889333,0.005348220947889334,3,1,50.59256466666667,151.777694,0.3572423751649123,0.3572423751649123,1.534523078904667,13.810707710142,0.001910611154330724,0.001910611154330724),
    @("ECs","Angptl3","Itgav","MuSCs",1,0.3333333333333333,0.030331,0.090993,0.005348220947889333,0.005348220947889334,3,1,26.84076266666667,80.522288,0.1895270158659356,0.1895270158659356,0.8141071724426667,7.326964551984001,0.001013632356445151,0.001013632356445151),
    @("ECs","Angptl3","Itgav","Resolving-Mac",1,0.3333333333333333,0.030331,0.090993,0.005348220947889333,0.005348220947889334,3,1,54.762539,164.287617,0.3866872460670236,0.3866872460670236,1.661002570409,14.949023133681,0.002068088829697293,0.002068088829697293),
    @("FAPs","Angptl3","Itgav","ECs",3,1,3.081251,9.243753,0.5433124903202978,0.5433124903202979,3,1,9.423852333333334,28.271557,0.06654336290212845,0.06654336290212845,29.03725442593567,261.335289833421,0.03615384021264273,0.03615384021264274),
    @("FAPs","Angptl3","Itgav","FAPs",3,1,3.081251,9.243753,0.5433124903202978,0.5433124903202979,3,1,50.59256466666667,151.777694,0.3572423751649123,0.3572423751649123,155.8883904717313,1402.995514245582,0.1940942444987866,0.1940942444987866),
    @("FAPs","Angptl3","Itgav","MuSCs",3,1,3.081251,9.243753,0.5433124903202978,0.5433124903202979,3,1,26.84076266666667,80.522288,0.1895270158659356,0.1895270158659356,82.70312680742933,744.328141266864,0.102972394973096,0.1029723949730961),
    @("FAPs","Angptl3","Itgav","Resolving-Mac",3,1,3.081251,9.243753,0.5433124903202978,0.5433124903202979,3,1,54.762539,164.287617,0.3866872460670236,0.3866872460670236,168.737128056289,1518.634152506601,0.2100920106357724,0.2100920106357725),
    @("MuSCs","Angptl3","Itgav","ECs",3,1,1.424064666666667,4.272194,0.2511032435928821,0.2511032435928821,3,1,9.423852333333334,28.271557,0.06654336290212845,0.06654336290212845,13.42017513178422,120.781576186058,0.01670925426430271,0.01670925426430271),
    @("MuSCs","Angptl3","Itgav","FAPs",3,1,1.424064666666667,4.272194,0.2511032435928821,0.2511032435928821,3,1,50.59256466666667,151.777694,0.3572423751649123,0.3572423751649123,72.04708373784844,648.4237536406359,0.08970471915273473,0.08970471915273476),
    @("MuSCs","Angptl3","Itgav","MuSCs",3,1,1.424064666666667,4.272194,0.2511032435928821,0.2511032435928821,3,1,26.84076266666667,80.522288,0.1895270158659356,0.1895270158659356,38.22298173998578,344.006835659872,0.04759084843241604,0.04759084843241606),
    @("MuSCs","Angptl3","Itgav","Resolving-Mac",3,1,1.424064666666667,4.272194,0.2511032435928821,0.2511032435928821,3,1,54.762539,164.287617,0.3866872460670236,0.3866872460670236,77.98539684685532,701.8685716216979,0.09709842174342856,0.09709842174342859),
    @("Resolving-Mac","Angptl3","Itgav","ECs",3,1,1.135585,3.406755,0.2002360451389308,0.2002360451389308,3,1,9.423852333333334,28.271557,0.06654336290212845,0.06654336290212845,10.70158535194834,96.31426816753502,0.01332437981776684,0.01332437981776684),
    @("Resolving-Mac","Angptl3","Itgav","FAPs",3,1,1.135585,3.406755,0.2002360451389308,0.2002360451389308,3,1,50.59256466666667,151.777694,0.3572423751649123,0.3572423751649123,57.45215754699667,517.06941792297,0.07153280035906022,0.07153280035906022),
    @("Resolving-Mac","Angptl3","Itgav","MuSCs",3,1,1.135585,3.406755,0.2002360451389308,0.2002360451389308,3,1,26.84076266666667,80.522288,0.1895270158659356,0.1895270158659356,30.47996747282667,274.3197072554401,0.03795014010397832,0.03795014010397833),
    @("Resolving-Mac","Angptl3","Itgav","Resolving-Mac",3,1,1.135585,3.406755,0.2002360451389308,0.2002360451389308,3,1,54.762539,164.287617,0.3866872460670236,0.3866872460670236,62.187517850315,559.687660652835,0.07742872485812537,0.07742872485812537)
)

$rowCount = $data.Count
$colCount = $data[0].Count

$arr = New-Object 'object[,]' $rowCount,$colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $arr[$r,$c] = $data[$r][$c]
    }
}

$startRow = 2
$endRow = $startRow + $rowCount - 1
$rng = $ws.Range($ws.Cells.Item($startRow,1), $ws.Cells.Item($endRow,20))
$rng.Value = $arr
